$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D2:D51) keeps its original text formatting so that
# numeric-looking values (e.g. "0.9998") are stored as text, not converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '28.986.51'
$ws.Range("E2").Value = '  +1.19%  '

# Row 3
$ws.Range("D3").Value = '1.888.33'
$ws.Range("E3").Value = '  +0.87%  '

# Row 4
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").Value = '331.24'
$ws.Range("E5").Value = '  -2.15%  '

# Row 6
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").Value = '0.4610'
$ws.Range("E7").Value = '  -1.87%  '

# Row 8
$ws.Range("D8").Value = '0.4100'
$ws.Range("E8").Value = '  +2.83%  '

# Row 9
$ws.Range("D9").Value = '47.37'
$ws.Range("E9").Value = '  -0.68%  '

# Row 10
$ws.Range("D10").Value = '0.07995'
$ws.Range("E10").Value = '  -0.75%  '

# Row 11
$ws.Range("D11").Value = '0.9907'
$ws.Range("E11").Value = '  -1.21%  '

# Row 12
$ws.Range("D12").Value = '21.73'
$ws.Range("E12").Value = '  -1.83%  '

# Row 13
$ws.Range("D13").Value = '1.876.11'
$ws.Range("E13").Value = '  +0.99%  '

# Row 14
$ws.Range("D14").Value = '5.910'
$ws.Range("E14").Value = '  -2.46%  '

# Row 15
$ws.Range("D15").Value = '7.076'
$ws.Range("E15").Value = '  -3.15%  '

# Row 16
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '88.99'
$ws.Range("E16").Value = '  -1.64%  '

# Row 17
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '0.9997'
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").Value = '0.00001030'
$ws.Range("E18").Value = '  -0.97%  '

# Row 19
$ws.Range("D19").Value = '0.06562'
$ws.Range("E19").Value = '  -0.98%  '

# Row 20
$ws.Range("D20").Value = '17.47'
$ws.Range("E20").Value = '  -0.91%  '

# Row 21
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.31%  '

# Row 22
$ws.Range("D22").Value = '29.030.22'
$ws.Range("E22").Value = '  +1.55%  '

# Row 23
$ws.Range("D23").Value = '5.406'
$ws.Range("E23").Value = '  -1.68%  '

# Row 24
$ws.Range("D24").Value = '11.25'
$ws.Range("E24").Value = '  +1.76%  '

# Row 25
$ws.Range("D25").Value = '2.208'
$ws.Range("E25").Value = '  -2.08%  '

# Row 26
$ws.Range("D26").Value = '2.133.54'
$ws.Range("E26").Value = '  +2.47%  '

# Row 27
$ws.Range("D27").Value = '157.28'
$ws.Range("E27").Value = '  -2.37%  '

# Row 28
$ws.Range("D28").Value = '19.64'
$ws.Range("E28").Value = '  -0.84%  '

# Row 29
$ws.Range("D29").Value = '2.109'

# Row 30
$ws.Range("D30").Value = '5.404'
$ws.Range("E30").Value = '  -1.65%  '

# Row 31
$ws.Range("D31").Value = '117.99'
$ws.Range("E31").Value = '  -1.56%  '

# Row 32
$ws.Range("D32").Value = '0.9771'
$ws.Range("E32").Value = '  +0.49%  '

# Row 33
$ws.Range("D33").Value = '0.09346'
$ws.Range("E33").Value = '  -2.02%  '

# Row 34
$ws.Range("D34").Value = '3.598'
$ws.Range("E34").Value = '  +0.28%  '

# Row 35
$ws.Range("D35").Value = '1.411'
$ws.Range("E35").Value = '  +1.79%  '

# Row 36
$ws.Range("D36").Value = '5.283'
$ws.Range("E36").Value = '  -1.79%  '

# Row 37
$ws.Range("D37").Value = '0.06050'
$ws.Range("E37").Value = '  -2.91%  '

# Row 38
$ws.Range("D38").Value = '0.02227'
$ws.Range("E38").Value = '  -1.11%  '

# Row 39
$ws.Range("D39").Value = '8.271'
$ws.Range("E39").Value = '  -1.19%  '

# Row 40
$ws.Range("D40").Value = '1.187'
$ws.Range("E40").Value = '  +0.27%  '

# Row 41
$ws.Range("D41").Value = '0.9993'
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").Value = '0.5776'
$ws.Range("E42").Value = '  -2.88%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '10.13'
$ws.Range("E43").Value = '  -2.13%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.1819'
$ws.Range("E44").Value = '  -3.34%  '

# Row 45
$ws.Range("D45").Value = '1.259'
$ws.Range("E45").Value = '  +0.01%  '

# Row 46
$ws.Range("D46").Value = '2.278'
$ws.Range("E46").Value = '  +9.55%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '12.01'
$ws.Range("E47").Value = '  -1.17%  '

# Row 48
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.5469'
$ws.Range("E48").Value = '  -1.83%  '

# Row 49
$ws.Range("D49").Value = '1.901'
$ws.Range("E49").Value = '  -3.05%  '

# Row 50
$ws.Range("D50").Value = '0.07020'
$ws.Range("E50").Value = '  -5.29%  '

# Row 51
$ws.Range("D51").Value = '45.52'
$ws.Range("E51").Value = '  +14.26%  '
